$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from column R (previous last year column) into column S for rows 3-14
$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the 2023 values
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 48.8
$ws.Range("S5").Value = 49.7
$ws.Range("S6").Value = 811
$ws.Range("S7").Value = 668
$ws.Range("S8").Value = 670.4
$ws.Range("S9").Value = 14.1
$ws.Range("S10").Value = 5.4
$ws.Range("S11").Value = 35.6
$ws.Range("S12").Value = 2.2
$ws.Range("S13").Value = 14.1
$ws.Range("S14").Value = 0

# Update the view: scroll so column I is leftmost, and select S4:S14 with active cell S4
$ws.Range("S4:S14").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 9  # column I
